$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2 through 451) holds the "Förändrad" date, stored as a
# serial date number. Every value moves from 45171 (2023-09-02) to
# 45172 (2023-09-03).
$ws.Range("C2:C451").Value = 45172
